# Model def file structure updates.
# Adds two new columns ("recovery_function" and "location") to the
# comp_type_dmg_algo sheet, fixes a stray top-border on the
# component_class (D) column so it matches its row group, and resets
# the active selection cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("comp_type_dmg_algo")

# ---------------------------------------------------------------------
# 1) Insert "recovery_function" column immediately before "recovery_mean"
#    (column M, 13th column) - done first so the new shared strings are
#    appended in the same order as the source edit (recovery_function,
#    then Normal).
# ---------------------------------------------------------------------
$ws.Columns.Item(13).Insert()
$ws.Range("M1").Value = "recovery_function"
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 13).Value = "Normal"
}

# ---------------------------------------------------------------------
# 2) Insert "location" column immediately before "beta" (column G, 7th
#    column).
# ---------------------------------------------------------------------
$ws.Columns.Item(7).Insert()
$ws.Range("G1").Value = "location"
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 7).Value = 0
}

# ---------------------------------------------------------------------
# 3) The component_class column (D) carried a stray top border on every
#    row; match it back up with the rest of its row group (column C).
# ---------------------------------------------------------------------
for ($r = 2; $r -le 25; $r++) {
    $src = $ws.Cells.Item($r, 3)
    $dst = $ws.Cells.Item($r, 4)
    $src.Copy()
    $dst.PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4) Reset the active cell/selection on this sheet.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("B1").Select()
